$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Write-Host "Before insert, dimension: $($ws.UsedRange.Address())"
$ws.Rows("180:181").Insert()
Write-Host "After insert, dimension: $($ws.UsedRange.Address())"
Write-Host "K180=$($ws.Cells.Item(180,11).Value())"
Write-Host "K182=$($ws.Cells.Item(182,11).Value())"
Write-Host "K184=$($ws.Cells.Item(184,11).Value())"
